$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new ImagePath column (bold, matching the other header cells)
$ws.Range("D1").Value = "ImagePath"
$ws.Range("D1").Font.Bold = $true

# Build the ImagePath value for each data row (2-46) based on the Item Name in column A
for ($r = 2; $r -le 46; $r++) {
    $itemName = $ws.Cells.Item($r, 1).Text
    $imageName = "tex_DaS_" + $itemName.Replace(" ", "") + ".png"
    $ws.Cells.Item($r, 4).Value = $imageName
}

# Match new column widths for the ImagePath column (D) and the now-wider adjacent column (E)
$ws.Columns.Item(4).ColumnWidth = 66.7109375
$ws.Columns.Item(5).ColumnWidth = 65.42578125

# Update the sheet's active selection
$ws.Range("E2:E46").Select()
